# Fix update progress template: insert a new "Activity Name" column
# between "Resource Code" (B) and "Progress" (C), pushing "Progress" to D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C (shifts the existing "Progress" column to D).
# The new column inherits the header style (s=3) from column B automatically.
$ws.Range("C1").EntireColumn.Insert()

# Set the new header's text.
$ws.Range("C1").Value = "Activity Name"

# Match the width of the surrounding header columns (A/B use 24 / 22.66).
$ws.Columns("C:C").ColumnWidth = 23.17
